$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell labels (renaming of field/link identifiers)
$ws.Range("D1").Value = "div_testingTools_internalText"
$ws.Range("E1").Value = "div_testingTools_internalText_1"
$ws.Range("F1").Value = "link_projectLinks_internalRoleLinkName"
$ws.Range("G1").Value = "link_projectLinks_project_id"
$ws.Range("H1").Value = "link_projectLinks_team_id"

# Adjust column widths for columns D:H.
# Excel's ColumnWidth property is offset from the stored OOXML column
# width by the default font's padding (~0.8333333333333334 chars for
# this workbook's default Calibri 11 font), so subtract that offset to
# land exactly on the target stored widths of 31, 33, 40, 30, 27.
$widthOffset = 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 31 - $widthOffset
$ws.Columns.Item(5).ColumnWidth = 33 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 40 - $widthOffset
$ws.Columns.Item(7).ColumnWidth = 30 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 27 - $widthOffset
